$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) First paragraph: "This is a Microsoft word document." gets two
#    trailing spaces, then a new red "(This is a change - Version for
#    main branch)" annotation appended as three separate runs (all
#    sharing the same red formatting, mirroring the authored edit,
#    which recorded them as three <w:r> elements).
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$p1TextRange = $d.Range($p1.Range.Start, $p1.Range.End - 1)

# Append two trailing spaces to the existing (default-formatted) run.
$p1TextRange.InsertAfter("  ")

# Recompute where the plain-text portion now ends (just before the
# paragraph mark) so the new runs get inserted right after it.
$p1 = $d.Paragraphs.Item(1)
$insertAt = $p1.Range.End - 1

$red = 255  # wdColorRed / RGB(255,0,0) -> OOXML w:val="FF0000"

$seg1 = "(This is a change " + [char]0x2013 + " Ve"
$r1 = $d.Range($insertAt, $insertAt)
$r1.InsertAfter($seg1)
$r1 = $d.Range($insertAt, $insertAt + $seg1.Length)
$r1.Font.Color = $red
$insertAt = $insertAt + $seg1.Length

$seg2 = "rsion for main branch"
$r2 = $d.Range($insertAt, $insertAt)
$r2.InsertAfter($seg2)
$r2 = $d.Range($insertAt, $insertAt + $seg2.Length)
$r2.Font.Color = $red
$insertAt = $insertAt + $seg2.Length

$seg3 = ")"
$r3 = $d.Range($insertAt, $insertAt)
$r3.InsertAfter($seg3)
$r3 = $d.Range($insertAt, $insertAt + $seg3.Length)
$r3.Font.Color = $red

# ---------------------------------------------------------------------
# 2) Drop the trailing "ank God almighty, we are free at last."
#    paragraph entirely (it followed "Shall be lifted-nevermore!" and
#    is removed in full, including its own paragraph mark).
# ---------------------------------------------------------------------
$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastIndex)
$lastPara.Range.Delete()

# ---------------------------------------------------------------------
# 3) Styles that are no longer referenced anywhere in the body are
#    dropped from the style sheet, mirroring the authored save. Delete
#    from the end of the collection backwards so earlier ordinals stay
#    valid while later ones are removed.
# ---------------------------------------------------------------------
$unusedStyleNames = @(
  "podcast-toolssubscribe-links",
  "generic-title",
  "subscribe-more-info",
  "subscribe",
  "audio-tool",
  "Heading4Char",
  "Heading2Char",
  "Hyperlink",
  "apple-converted-space",
  "Heading4",
  "Heading2"
)
foreach ($styleName in $unusedStyleNames) {
  $style = $d.Styles.Item($styleName)
  if ($style -ne $null) {
    $style.Delete()
  }
}
